# Generate Report for Handoff
# Rotates the localization-status report onto the next handback run:
#   - new source markdown guid (4db14c90... -> 759b30f8...)
#   - new xliff file names (new guid + new content hash)
#   - refreshed handoff/generate timestamps
#   - the zh-cn / de-de sheets no longer have a "Latest Target File" handed
#     back yet, so those columns (and their hyperlink) are cleared

$wb = $excel.ActiveWorkbook

$oldGuid = "4db14c90-50ac-469d-8c12-c56368b6f730"
$newGuid = "759b30f8-4682-4803-baaf-463f9a13eaac"

$newMdName = "$newGuid.md"
$newMdPath = "e2e\$newGuid.md"

$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1c1397601d9442fc55afa3b5232f8d0ef3afdba/e2e/$oldGuid.md"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("G2").Value = "2016-09-07 03:16:00"

# Refresh the B2 hyperlink display text, keeping the same link target.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $mdUrl, "", "", $newMdPath) | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newMdName
$wsZhCn.Range("G2").Value = "$newGuid.fbf062b63eade5cd6321f9fb68a2227fc8f725ac.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-07 03:15:54"
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

# The handback file/name for this run hasn't been produced yet -> clear it,
# and drop its hyperlink along with the old "Latest Target File" value.
$wsZhCn.Range("I2").Style = "Normal"
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("J2").Value = ""

# Re-create the remaining A2 hyperlink (same target, new display text) after
# clearing the sheet's hyperlinks collection above.
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdUrl, "", "", $newMdName) | Out-Null

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newMdName
$wsDeDe.Range("G2").Value = "$newGuid.fbf062b63eade5cd6321f9fb68a2227fc8f725ac.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-07 03:16:00"
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDeDe.Range("I2").Style = "Normal"
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("J2").Value = ""

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdUrl, "", "", $newMdName) | Out-Null

# ---------------------------------------------------------------------
# Column width tweaks on zh-cn / de-de: columns I/J no longer need to be
# wide enough for the guid-based handback filenames.
# ---------------------------------------------------------------------
foreach ($ws in @($wsZhCn, $wsDeDe)) {
    $ws.Range("I1").ColumnWidth = 18.6506053379604
    $ws.Range("J1").ColumnWidth = 21.7054770333426
}
